# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, reflecting refreshed scrape output.
# Generated to match commit: "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 3697
$ws1.Range("F5").Value  = 3697
$ws1.Range("F7").Value  = 5225
$ws1.Range("F8").Value  = 563
$ws1.Range("F9").Value  = 395
$ws1.Range("F11").Value = 722
$ws1.Range("F16").Value = 337
$ws1.Range("F17").Value = 42
$ws1.Range("F21").Value = 367
$ws1.Range("F22").Value = 5979
$ws1.Range("F26").Value = 6292
$ws1.Range("F29").Value = 3241
$ws1.Range("F30").Value = 358
$ws1.Range("F31").Value = 735
$ws1.Range("F33").Value = 322
$ws1.Range("F35").Value = 147
$ws1.Range("F36").Value = 1097
$ws1.Range("F37").Value = 94
$ws1.Range("F40").Value = 905
$ws1.Range("F41").Value = 1076
$ws1.Range("F42").Value = 2045

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 29

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1141

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1141
$ws4.Range("F7").Value  = 3697
$ws4.Range("F8").Value  = 3697
$ws4.Range("F10").Value = 5225
$ws4.Range("F11").Value = 563
$ws4.Range("F12").Value = 395
$ws4.Range("F14").Value = 722
$ws4.Range("F19").Value = 337
$ws4.Range("F20").Value = 42
$ws4.Range("F25").Value = 367
$ws4.Range("F26").Value = 5979
$ws4.Range("F30").Value = 6292
$ws4.Range("F33").Value = 3241
$ws4.Range("F34").Value = 358
$ws4.Range("F35").Value = 735
$ws4.Range("F37").Value = 322
$ws4.Range("F38").Value = 29
$ws4.Range("F40").Value = 147
$ws4.Range("F41").Value = 1097
$ws4.Range("F42").Value = 94
$ws4.Range("F45").Value = 905
$ws4.Range("F46").Value = 1076
$ws4.Range("F48").Value = 2045
